# Update "想去人数" (interest/attendee counts) for two sheets: "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F6").Value = 38
    $ws.Range("F7").Value = 121
    $ws.Range("F9").Value = 299
}
